$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.240.51"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.272.30"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.78"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.34"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.268.39"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.47"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.20"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.818.35"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.276.57"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.370.41"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.79"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.85"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.68"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.19"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.07"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.32"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.49"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.05"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0402"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "427.36"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.061.29"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("E42").Value = "  +6.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.22"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.34"
$ws.Range("E47").Value = "  +7.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.10"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.81"
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("E51").Value = "  -1.08%  "
